{"js": "// Update the date line and the 25 \"two-digit \u00f7 one-digit\" problems in the\n// table, replacing each value in place while keeping existing formatting.\n\n// 1) Date paragraph (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst datePara = paragraphs.items[0];\ndatePara.load(\"text\");\nawait context.sync();\n\nif (datePara.text === \"2025-10-21 Tuesday\") {\n  datePara.insertText(\"2025-10-22 Wednesday\", Word.InsertLocation.replace);\n}\n\n// 2) Table cells, in row-major reading order; only the cells that actually\n// hold a problem (\"N\u00f7N=\") need to change \u2014 the new value replaces the old\n// one positionally (several problems repeat their original text, e.g.\n// \"57\u00f74=\" and \"32\u00f78=\" each occur twice with different replacements, so this\n// must be driven by position, not by searching for the old text).\nconst replacements = [\n  \"78\u00f75=\", \"32\u00f78=\", \"68\u00f78=\", \"91\u00f72=\", \"80\u00f75=\",\n  \"50\u00f75=\", \"80\u00f74=\", \"29\u00f79=\", \"72\u00f74=\", \"13\u00f77=\",\n  \"49\u00f76=\", \"18\u00f73=\", \"12\u00f73=\", \"95\u00f79=\", \"54\u00f75=\",\n  \"97\u00f72=\", \"36\u00f75=\", \"96\u00f72=\", \"17\u00f73=\", \"24\u00f79=\",\n  \"22\u00f76=\", \"25\u00f73=\", \"61\u00f72=\", \"89\u00f75=\", \"67\u00f74=\"\n];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"isNullObject\");\nawait context.sync();\n\nif (!table.isNullObject) {\n  const rows = table.rows;\n  rows.load(\"items\");\n  await context.sync();\n\n  for (const row of rows.items) {\n    row.cells.load(\"items\");\n  }\n  await context.sync();\n\n  // Use each cell's first paragraph (not cell.body) so insertText keeps the\n  // existing paragraph/run formatting (alignment, font, size) instead of\n  // resetting it.\n  const cellParagraphs = [];\n  for (const row of rows.items) {\n    for (const cell of row.cells.items) {\n      cellParagraphs.push(cell.body.paragraphs.getFirst());\n    }\n  }\n  for (const para of cellParagraphs) {\n    para.load(\"text\");\n  }\n  await context.sync();\n\n  let replIndex = 0;\n  for (const para of cellParagraphs) {\n    const text = para.text.trim();\n    if (text.length > 0 && replIndex < replacements.length) {\n      para.insertText(replacements[replIndex], Word.InsertLocation.replace);\n      replIndex++;\n    }\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the 25 \"two-digit \u00f7 one-digit\" problems in the\n# table, replacing each value in place while keeping existing formatting.\n\n$d = $word.ActiveDocument\n\n# 1) Date paragraph (first paragraph in the document).\n$dateRange = $d.Paragraphs(1).Range\nif ($dateRange.Text.TrimEnd([char]13, [char]7) -eq \"2025-10-21 Tuesday\") {\n    $dateRange.Text = \"2025-10-22 Wednesday\"\n}\n\n# 2) Table cells, in row-major reading order; only the 25 cells that hold a\n# problem (\"N\u00f7N=\") need to change. Several original values repeat (e.g.\n# \"57\u00f74=\" and \"32\u00f78=\" each occur twice with different replacements), so the\n# mapping must be driven by position, not by searching for the old text.\n$replacements = @(\n    \"78\u00f75=\", \"32\u00f78=\", \"68\u00f78=\", \"91\u00f72=\", \"80\u00f75=\",\n    \"50\u00f75=\", \"80\u00f74=\", \"29\u00f79=\", \"72\u00f74=\", \"13\u00f77=\",\n    \"49\u00f76=\", \"18\u00f73=\", \"12\u00f73=\", \"95\u00f79=\", \"54\u00f75=\",\n    \"97\u00f72=\", \"36\u00f75=\", \"96\u00f72=\", \"17\u00f73=\", \"24\u00f79=\",\n    \"22\u00f76=\", \"25\u00f73=\", \"61\u00f72=\", \"89\u00f75=\", \"67\u00f74=\"\n)\n\n$table = $d.Tables(1)\n$rowCount = $table.Rows.Count\n$colCount = $table.Columns.Count\n\n$i = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cellRange = $table.Cell($r, $c).Range\n        $text = $cellRange.Text.TrimEnd([char]13, [char]7)\n        if ($text.Length -gt 0 -and $i -lt $replacements.Length) {\n            $cellRange.Text = $replacements[$i]\n            $i++\n        }\n    }\n}\n"}
